# Apply minor numeric recalculation updates to the COD_summary sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J2").Value = 13.3808623608729

$ws.Range("I3").Value = 1774.875622030926
$ws.Range("J3").Value = 12.74777477556336
$ws.Range("K3").Value = 1.159414701257576
$ws.Range("N3").Value = 121.380862360873

$ws.Range("I4").Value = 1774.607478521407
$ws.Range("J4").Value = 12.45488223305614
$ws.Range("K4").Value = 1.233763307075057
$ws.Range("L4").Value = 2.131015916978649
$ws.Range("M4").Value = 0.06556972052241997
$ws.Range("N4").Value = 134.1286371364363
$ws.Range("O4").Value = 56.23941470125762

$ws.Range("I5").Value = 1774.3394191673
$ws.Range("J5").Value = 12.20629766163735
$ws.Range("K5").Value = 1.30528449468716
$ws.Range("L5").Value = 2.042023219613001
$ws.Range("M5").Value = 0.06283148368040004
$ws.Range("N5").Value = 146.5835193694924
$ws.Range("O5").Value = 57.47317800833267

$ws.Range("I6").Value = 1774.071437309771
$ws.Range("J6").Value = 11.99866664560645
$ws.Range("K6").Value = 1.374320712981774
$ws.Range("L6").Value = 3.156425100544646
$ws.Range("M6").Value = 0.09712077232445063
$ws.Range("N6").Value = 158.7898170311298
$ws.Range("O6").Value = 58.77846250301983

$ws.Range("I7").Value = 1773.803526794177
$ws.Range("J7").Value = 11.47324891637379
$ws.Range("K7").Value = 1.494693836948371
$ws.Range("L7").Value = 1.890639607614224
$ws.Range("M7").Value = 0.05817352638812998
$ws.Range("N7").Value = 170.7884836767364
$ws.Range("O7").Value = 60.1527832160016

$ws.Range("I8").Value = 1773.535735277823
$ws.Range("J8").Value = 11.35245243162619
$ws.Range("K8").Value = 1.566858109687555
$ws.Range("L8").Value = 2.482194601935355
$ws.Range("M8").Value = 0.07637521852108783
$ws.Range("N8").Value = 182.2617325931104
$ws.Range("O8").Value = 61.64747705294996

$ws.Range("I9").Value = 1773.268002049666
$ws.Range("J9").Value = 11.09650387758536
$ws.Range("K9").Value = 1.659139473800495
$ws.Range("L9").Value = 2.172733351266169
$ws.Range("M9").Value = 0.0668533338851129
$ws.Range("N9").Value = 193.6141850247366
$ws.Range("O9").Value = 63.2143351626375

$ws.Range("I10").Value = 1773.000347373775
$ws.Range("J10").Value = 10.94959980496654
$ws.Range("K10").Value = 1.780116714042583
$ws.Range("L10").Value = 1.540244479045344
$ws.Range("M10").Value = 0.04739213781677981
$ws.Range("N10").Value = 204.7106889023218
$ws.Range("O10").Value = 64.87347463643796

$ws.Range("C11").Value = 11818218
$ws.Range("I11").Value = 1772.732754881697
$ws.Range("J11").Value = 10.98443980637041
$ws.Range("K11").Value = 1.853240598998291
$ws.Range("L11").Value = 2.569908419474952
$ws.Range("M11").Value = 0.0790741052146139
$ws.Range("N11").Value = 215.6602887072884
$ws.Range("O11").Value = 66.65359135048054
